$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "16÷7="
$tbl.Cell(1,2).Range.Text = "92÷5="
$tbl.Cell(1,3).Range.Text = "67÷6="
$tbl.Cell(1,4).Range.Text = "36÷9="
$tbl.Cell(1,5).Range.Text = "52÷2="
$tbl.Cell(5,1).Range.Text = "13÷7="
$tbl.Cell(5,2).Range.Text = "50÷6="
$tbl.Cell(5,3).Range.Text = "46÷9="
$tbl.Cell(5,4).Range.Text = "80÷7="
$tbl.Cell(5,5).Range.Text = "23÷9="
$tbl.Cell(9,1).Range.Text = "60÷5="
$tbl.Cell(9,2).Range.Text = "89÷3="
$tbl.Cell(9,3).Range.Text = "33÷3="
$tbl.Cell(9,4).Range.Text = "53÷3="
$tbl.Cell(9,5).Range.Text = "55÷3="
$tbl.Cell(13,1).Range.Text = "94÷9="
$tbl.Cell(13,2).Range.Text = "13÷3="
$tbl.Cell(13,3).Range.Text = "36÷9="
$tbl.Cell(13,4).Range.Text = "75÷8="
$tbl.Cell(13,5).Range.Text = "50÷6="
$tbl.Cell(17,1).Range.Text = "18÷6="
$tbl.Cell(17,2).Range.Text = "85÷3="
$tbl.Cell(17,3).Range.Text = "50÷7="
$tbl.Cell(17,4).Range.Text = "67÷8="
$tbl.Cell(17,5).Range.Text = "40÷6="
